$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A100").Value = "i,e,f,b,g"
$ws.Range("A101").Value = "i,e,f,b,g"
$ws.Range("A102").Value = "c,i,i,i,i,i,f,h,b,g"
$ws.Range("A103").Value = "c,i,i,i,i,i,f,h,b,g"

$ws.Range("A126").Select()
